$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-12-06"

# Update the row label for December
$ws.Range("A14").Value = "December (through 12-06)"

# Row 14 updates (2016 + 2020 + totals columns)
$ws.Range("F14").Value = 18
$ws.Range("G14").Value = 0.0526
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 20
$ws.Range("J14").Value = 0.1304
$ws.Range("O14").Value = 6
$ws.Range("R14").Value = 29
$ws.Range("S14").Value = 0.0645
$ws.Range("U14").Value = 49

# Row 15 (Total) updates
$ws.Range("F15").Value = 521
$ws.Range("G15").Value = 0.1048
$ws.Range("H15").Value = 66
$ws.Range("I15").Value = 778
$ws.Range("J15").Value = 0.0782
$ws.Range("O15").Value = 486
$ws.Range("P15").Value = 0.1
$ws.Range("R15").Value = 1229
$ws.Range("S15").Value = 0.051
$ws.Range("U15").Value = 1593
$ws.Range("V15").Value = 0.0585
